# Minor update to documentation (slide 6 - "Lore - Next Steps")
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Move the body placeholder up (y: 1089050 -> 944725 EMU) ---
$shp.Top = 74.3878

# --- Tighten line spacing on the first two paragraphs (115% -> 100%) ---
$tr.Paragraphs(1, 1).ParagraphFormat.SpaceWithin = 1.0
$tr.Paragraphs(2, 1).ParagraphFormat.SpaceWithin = 1.0

# --- Paragraph 2 (Database bullet): fix "the data geodata" -> "the geodata"
#     and "comparison f data" -> "comparison of data" (second run only) ---
$old2 = "Lore should include the option of using a database instead of tracking events in memory to prevent it from being overwhelmed by large datasets. While a NoSQL approach would work, the data geodata tends be highly structured, making a traditional SQL database an equally viable option. This would allow easy saving and comparison f data from multiple runs."
$new2 = "Lore should include the option of using a database instead of tracking events in memory to prevent it from being overwhelmed by large datasets. While a NoSQL approach would work, the geodata tends be highly structured, making a traditional SQL database an equally viable option. This would allow easy saving and comparison of data from multiple runs."
$full = $tr.Text
$idx2 = $full.IndexOf($old2)
$tr.Characters($idx2 + 1, $old2.Length).Text = $new2

# --- Paragraph 3: "Entity information" -> "Event information" ---
$para3 = $tr.Paragraphs(3, 1)
$tr.Characters($para3.Start, $para3.Length).Text = "Convert time from sim to real-world - Event information should display real world time, rather than sim time by default."

# --- Paragraph 5: append clarifying clause about crossfilter ---
$para5 = $tr.Paragraphs(5, 1)
$tr.Characters($para5.Start, $para5.Length).Text = "Add crossfilter - It should be possible to filter based on entity state, type, or process (allow filtering on subsets of the global statistics)."

# --- Paragraph 6: add trailing period to "... to a kml movie" run ---
$full = $tr.Text
$old6 = " to a kml movie"
$idx6 = $full.IndexOf($old6)
$tr.Characters($idx6 + 1, $old6.Length).Text = " to a kml movie."
